# Update a handful of numeric values (column F) across the "展览" and
# "全部类型" worksheets, and one value on "演出", per the upstream data
# refresh (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 1113
$wsExpo.Range("F8").Value = 397
$wsExpo.Range("F14").Value = 12728
$wsExpo.Range("F16").Value = 5240

# 演出 (sheet2)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 95

# 全部类型 (sheet4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1113
$wsAll.Range("F9").Value = 397
$wsAll.Range("F15").Value = 12728
$wsAll.Range("F16").Value = 95
$wsAll.Range("F19").Value = 5240
